$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Sponsored`nLaunching 2 BHK in Thane West - Pay 20% & Nothing till Jan'25`nraymondtenxera.com`nhttps://www.raymondtenxera.com › thane › project`nNew Launch Homes by Raymond with 38 Habitable Floors, 26500 SqFt Clubhouse, 40+ Amenities. Experience a futuristic lifestyle with Raymond Realty's Spacious 2 BHK homes..."
$ws.Range("B2").Value = "Raymond Limited"

$ws.Range("A3").Value = "Sponsored`n2 BHK Projects in Thane West | Starts at ₹93 Lacs* by Runwal`nlandsend.runwal.com`nhttp://landsend.runwal.com › projects › thane`nTake Advantage of the Umbrella Offer: 2 BHK Flats Starts at ₹93L* at Lands End by Runwal"

$ws.Range("A4").Value = "Sponsored`n1 BHK by Lodha® in Thane | 1,2,3 BHK by Lodha® in Thane`nLodha Group`nhttps://www.lodhagroup.in`nWorld-class 1 BHKs in Thane by India's #1 real estate developer. Building a better life."
$ws.Range("B4").Value = "Macrotech Developers Limited"

$ws.Range("A5").Value = "Sponsored`nGodrej Ascend, Kolshet, Thane | 2&3 BHK at ₹1.09Cr+*(All Incl)`nGodrej Properties`nhttps://www.godrejproperties.com`nReach Airoli in 20 min* and avail easy connectivity to Navi-Mumbai via Thane-Belapur Road"

$ws.Rows("6").Delete()
